$wb = $excel.ActiveWorkbook

# Row data to append (rows 144-146) for each of the 4 worksheets.
# Columns: A=timestamp(date), B/C/D/E=hex-string text, F/G/H/I=numbers

$rowsData = @{
    1 = @(
        @{ A = 45930.43784722222; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x00,0xF8"; E = "0x14"; F = 380; G = 759863127514710945038336.0; H = 248; I = 14 },
        @{ A = 45931.43920138889; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x00,0xF4"; E = "0x14"; F = 380; G = 759863127514710945038336.0; H = 244; I = 14 },
        @{ A = 45932.43804398148; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x00,0xF0"; E = "0x14"; F = 380; G = 759863127514710945038336.0; H = 240; I = 14 }
    )
    2 = @(
        @{ A = 45930.43784722222; B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x00,0xF8"; E = "0xe"; F = 380; G = 568432987514711010443264.0; H = 248; I = 14 },
        @{ A = 45931.43920138889; B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x00,0xF8"; E = "0xe"; F = 380; G = 568432987514711010443264.0; H = 248; I = 14 },
        @{ A = 45932.43804398148; B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x00,0xF8"; E = "0xe"; F = 380; G = 568432987514711010443264.0; H = 248; I = 14 }
    )
    3 = @(
        @{ A = 45930.43784722222; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x6C"; E = "0x7"; F = 130; G = 568631262647113970876416.0; H = 108; I = 7 },
        @{ A = 45931.43920138889; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x6C"; E = "0x7"; F = 130; G = 568631262647113970876416.0; H = 108; I = 7 },
        @{ A = 45932.43804398148; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x6C"; E = "0x7"; F = 130; G = 568631262647113970876416.0; H = 108; I = 7 }
    )
    4 = @(
        @{ A = 45930.43784722222; B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x6B"; E = "0x3"; F = 130; G = 985046333984776009023488.0; H = 107; I = 3 },
        @{ A = 45931.43920138889; B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x6B"; E = "0x3"; F = 130; G = 985046333984776009023488.0; H = 107; I = 3 },
        @{ A = 45932.43804398148; B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x6A"; E = "0x3"; F = 130; G = 985046333984776009023488.0; H = 106; I = 3 }
    )
}

for ($sheetIdx = 1; $sheetIdx -le 4; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $startRow = 144
    $rows = $rowsData[$sheetIdx]

    for ($k = 0; $k -lt $rows.Count; $k++) {
        $r = $startRow + $k
        $d = $rows[$k]

        $ws.Cells.Item($r, 1).Value = $d.A
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

        $ws.Cells.Item($r, 2).Value = $d.B
        $ws.Cells.Item($r, 3).Value = $d.C
        $ws.Cells.Item($r, 4).Value = $d.D
        $ws.Cells.Item($r, 5).Value = $d.E

        $ws.Cells.Item($r, 6).Value = $d.F
        $ws.Cells.Item($r, 7).Value = $d.G
        $ws.Cells.Item($r, 8).Value = $d.H
        $ws.Cells.Item($r, 9).Value = $d.I
    }
}
